$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = "<吴岳恒>yueheng.wu@net263.com,<李浩>hao1.li@net263.com,<KeithXi>xiaofeng.xi@net263.com"

$wb.Save()
